# Weekly update: insert a new data row for "Haba" (Feria Lagunitas de Puerto Montt)
# at row 36, pushing the existing rows 36-70 down to 37-71.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36 (shifts old rows 36..70 down to 37..71,
# carrying their data, formulas and formatting with them).
$ws.Rows(36).Insert()

# Populate the newly inserted row 36 with this week's new record.
$ws.Cells.Item(36, 1).Value = 4
$ws.Cells.Item(36, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(36, 3).Value = "Los Lagos"
$ws.Cells.Item(36, 4).Value = 44539
$ws.Cells.Item(36, 5).Value = 10
$ws.Cells.Item(36, 6).Value = 100112026
$ws.Cells.Item(36, 7).Value = "Haba"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 80
$ws.Cells.Item(36, 11).Value = 12000
$ws.Cells.Item(36, 12).Value = 12000
$ws.Cells.Item(36, 13).Value = 12000
$ws.Cells.Item(36, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(36, 16).Value = 480
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
